$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 0  # H40: 800 -> 0
$ws.Cells.Item(40, 9).Value = 0  # I40: 800 -> 0
$ws.Cells.Item(40, 11).Value = 0  # K40: 800 -> 0
$ws.Cells.Item(40, 13).ClearContents()  # M40 was -625
$ws.Cells.Item(98, 8).Value = 1918.5  # H98: 1863.421 -> 1918.5
$ws.Cells.Item(98, 9).Value = 1548.9286  # I98: 1503.8 -> 1548.9286
$ws.Cells.Item(98, 11).Value = 1548.9286  # K98: 1503.8 -> 1548.9286
$ws.Cells.Item(98, 13).Value = -50.92859999999996  # M98: -5.799999999999955 -> -50.92859999999996
$ws.Cells.Item(100, 8).Value = 2870.2307  # H100: 2871.923 -> 2870.2307
$ws.Cells.Item(100, 9).Value = 647.8570999999999  # I100: 651 -> 647.8570999999999
$ws.Cells.Item(100, 11).Value = 647.8570999999999  # K100: 651 -> 647.8570999999999
$ws.Cells.Item(100, 13).Value = -106.8570999999999  # M100: -110 -> -106.8570999999999
$ws.Cells.Item(107, 8).Value = 1514.1666  # H107: 712.6923 -> 1514.1666
$ws.Cells.Item(107, 9).Value = 1128  # I107: 392.66666 -> 1128
$ws.Cells.Item(107, 10).Value = 1900.3334  # J107: 1432.75 -> 1900.3334
$ws.Cells.Item(107, 11).Value = 1128  # K107: 392.66666 -> 1128
$ws.Cells.Item(107, 12).Value = 1900.3334  # L107: 1432.75 -> 1900.3334
$ws.Cells.Item(107, 13).Value = 792  # M107: 1527.33334 -> 792
$ws.Cells.Item(107, 14).Value = -5740.3334  # N107: -5272.75 -> -5740.3334
$ws.Cells.Item(116, 8).Value = 5699.5  # H116: 5731.6665 -> 5699.5
$ws.Cells.Item(116, 9).Value = 5699.5  # I116: 5731.6665 -> 5699.5
$ws.Cells.Item(116, 11).Value = 5699.5  # K116: 5731.6665 -> 5699.5
$ws.Cells.Item(116, 13).Value = -2257.5  # M116: -2289.6665 -> -2257.5
$ws.Cells.Item(122, 8).Value = 1918.5  # H122: 1863.421 -> 1918.5
$ws.Cells.Item(122, 9).Value = 1548.9286  # I122: 1503.8 -> 1548.9286
$ws.Cells.Item(122, 11).Value = 4646.7858  # K122: 4511.4 -> 4646.7858
$ws.Cells.Item(122, 13).Value = -2196.7858  # M122: -2061.4 -> -2196.7858
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 13497.333  # H17: 18896.8 -> 13497.333
$ws.Cells.Item(17, 9).Value = 3000  # I17: 0 -> 3000
$ws.Cells.Item(17, 10).Value = 15596.8  # J17: 18896.8 -> 15596.8
$ws.Cells.Item(17, 11).Value = 3000  # K17: 0 -> 3000
$ws.Cells.Item(17, 12).Value = 15596.8  # L17: 18896.8 -> 15596.8
$ws.Cells.Item(17, 13).Value = -2826  # M17: None -> -2826
$ws.Cells.Item(17, 14).Value = -15944.8  # N17: -19244.8 -> -15944.8
$ws.Cells.Item(134, 8).Value = 2273.238  # H134: 2339.8096 -> 2273.238
$ws.Cells.Item(134, 9).Value = 2494.5  # I134: 2563.611 -> 2494.5
$ws.Cells.Item(134, 10).Value = 945.6667  # J134: 997 -> 945.6667
$ws.Cells.Item(134, 11).Value = 7483.5  # K134: 7690.833 -> 7483.5
$ws.Cells.Item(134, 12).Value = 2837.0001  # L134: 2991 -> 2837.0001
$ws.Cells.Item(134, 13).Value = -4948.5  # M134: -5155.833 -> -4948.5
$ws.Cells.Item(134, 14).Value = -7907.0001  # N134: -8061 -> -7907.0001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 640.1429000000001  # H2: 560.375 -> 640.1429000000001
$ws.Cells.Item(2, 9).Value = 86.333336  # I2: 65.25 -> 86.333336
$ws.Cells.Item(2, 11).Value = 518.000016  # K2: 391.5 -> 518.000016
$ws.Cells.Item(2, 13).Value = -405.000016  # M2: -278.5 -> -405.000016
$ws.Cells.Item(12, 8).Value = 48.6  # H12: 52.555557 -> 48.6
$ws.Cells.Item(12, 10).Value = 46.57143  # J12: 52.166668 -> 46.57143
$ws.Cells.Item(12, 12).Value = 139.71429  # L12: 156.500004 -> 139.71429
$ws.Cells.Item(12, 14).Value = -485.71429  # N12: -502.500004 -> -485.71429
$ws.Cells.Item(44, 8).Value = 1399.6  # H44: 0 -> 1399.6
$ws.Cells.Item(44, 9).Value = 1499.5  # I44: 0 -> 1499.5
$ws.Cells.Item(44, 10).Value = 1000  # J44: 0 -> 1000
$ws.Cells.Item(44, 11).Value = 4498.5  # K44: 0 -> 4498.5
$ws.Cells.Item(44, 12).Value = 3000  # L44: 0 -> 3000
$ws.Cells.Item(44, 13).Value = -4100.5  # M44: None -> -4100.5
$ws.Cells.Item(44, 14).Value = -3796  # N44: None -> -3796
$ws.Cells.Item(131, 8).Value = 1311.875  # H131: 1427.8572 -> 1311.875
$ws.Cells.Item(131, 9).Value = 899.6  # I131: 999.5 -> 899.6
$ws.Cells.Item(131, 11).Value = 2698.8  # K131: 2998.5 -> 2698.8
$ws.Cells.Item(131, 13).Value = 2341.2  # M131: 2041.5 -> 2341.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 4500  # H113: 5000 -> 4500
$ws.Cells.Item(113, 9).Value = 4500  # I113: 5000 -> 4500
$ws.Cells.Item(113, 11).Value = 4500  # K113: 5000 -> 4500
$ws.Cells.Item(113, 13).Value = -2330  # M113: -2830 -> -2330
$ws.Cells.Item(122, 8).Value = 8931589  # H122: 8931837 -> 8931589
$ws.Cells.Item(122, 9).Value = 11365967  # I122: 12502862 -> 11365967
$ws.Cells.Item(122, 10).Value = 5533.3335  # J122: 4274.75 -> 5533.3335
$ws.Cells.Item(122, 11).Value = 34097901  # K122: 37508586 -> 34097901
$ws.Cells.Item(122, 12).Value = 16600.0005  # L122: 12824.25 -> 16600.0005
$ws.Cells.Item(122, 13).Value = -34095451  # M122: -37506136 -> -34095451
$ws.Cells.Item(122, 14).Value = -21500.0005  # N122: -17724.25 -> -21500.0005
$ws.Cells.Item(7, 8).Value = 3947  # H7: 0 -> 3947
$ws.Cells.Item(7, 9).Value = 3904  # I7: 0 -> 3904
$ws.Cells.Item(7, 10).Value = 3990  # J7: 0 -> 3990
$ws.Cells.Item(7, 11).Value = 3904  # K7: 0 -> 3904
$ws.Cells.Item(7, 12).Value = 3990  # L7: 0 -> 3990
$ws.Cells.Item(7, 13).Value = -3792  # M7: None -> -3792
$ws.Cells.Item(7, 14).Value = -4214  # N7: None -> -4214
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1393.3  # H16: 1341.1428 -> 1393.3
$ws.Cells.Item(16, 9).Value = 1103.6666  # I16: 1341.1428 -> 1103.6666
$ws.Cells.Item(16, 10).Value = 4000  # J16: 0 -> 4000
$ws.Cells.Item(16, 11).Value = 1103.6666  # K16: 1341.1428 -> 1103.6666
$ws.Cells.Item(16, 12).Value = 4000  # L16: 0 -> 4000
$ws.Cells.Item(16, 13).Value = -933.6666  # M16: -1171.1428 -> -933.6666
$ws.Cells.Item(16, 14).Value = -4340  # N16: None -> -4340
$ws.Cells.Item(18, 8).Value = 10000  # H18: 2000 -> 10000
$ws.Cells.Item(18, 10).Value = 10000  # J18: 2000 -> 10000
$ws.Cells.Item(18, 12).Value = 10000  # L18: 2000 -> 10000
$ws.Cells.Item(18, 14).Value = -10344  # N18: -2344 -> -10344
$ws.Cells.Item(20, 8).Value = 1000  # H20: 0 -> 1000
$ws.Cells.Item(20, 10).Value = 1000  # J20: 0 -> 1000
$ws.Cells.Item(20, 12).Value = 1000  # L20: 0 -> 1000
$ws.Cells.Item(20, 14).Value = -1452  # N20: None -> -1452
$ws.Cells.Item(22, 8).Value = 3497.25  # H22: 4000 -> 3497.25
$ws.Cells.Item(22, 9).Value = 3000  # I22: 0 -> 3000
$ws.Cells.Item(22, 10).Value = 3663  # J22: 4000 -> 3663
$ws.Cells.Item(22, 11).Value = 3000  # K22: 0 -> 3000
$ws.Cells.Item(22, 12).Value = 3663  # L22: 4000 -> 3663
$ws.Cells.Item(22, 13).Value = -2705  # M22: None -> -2705
$ws.Cells.Item(22, 14).Value = -4253  # N22: -4590 -> -4253
$ws.Cells.Item(27, 8).Value = 3497.25  # H27: 4000 -> 3497.25
$ws.Cells.Item(27, 9).Value = 3000  # I27: 0 -> 3000
$ws.Cells.Item(27, 10).Value = 3663  # J27: 4000 -> 3663
$ws.Cells.Item(27, 11).Value = 3000  # K27: 0 -> 3000
$ws.Cells.Item(27, 12).Value = 3663  # L27: 4000 -> 3663
$ws.Cells.Item(27, 13).Value = -2893  # M27: None -> -2893
$ws.Cells.Item(27, 14).Value = -3877  # N27: -4214 -> -3877
$ws.Cells.Item(40, 8).Value = 8371.875  # H40: 8853.143 -> 8371.875
$ws.Cells.Item(40, 9).Value = 5658.6665  # I40: 5986.5 -> 5658.6665
$ws.Cells.Item(40, 11).Value = 5658.6665  # K40: 5986.5 -> 5658.6665
$ws.Cells.Item(40, 13).Value = -5522.6665  # M40: -5850.5 -> -5522.6665
$ws.Cells.Item(55, 8).Value = 156.25  # H55: 143.18182 -> 156.25
$ws.Cells.Item(55, 10).Value = 450  # J55: 600 -> 450
$ws.Cells.Item(55, 12).Value = 450  # L55: 600 -> 450
$ws.Cells.Item(55, 14).Value = -796  # N55: -946 -> -796
$ws.Cells.Item(61, 8).Value = 7418.8335  # H61: 7423 -> 7418.8335
$ws.Cells.Item(61, 9).Value = 7422.8  # I61: 7423 -> 7422.8
$ws.Cells.Item(61, 10).Value = 7399  # J61: 0 -> 7399
$ws.Cells.Item(61, 11).Value = 7422.8  # K61: 7423 -> 7422.8
$ws.Cells.Item(61, 12).Value = 7399  # L61: 0 -> 7399
$ws.Cells.Item(61, 13).Value = -7220.8  # M61: -7221 -> -7220.8
$ws.Cells.Item(61, 14).Value = -7803  # N61: None -> -7803
$ws.Cells.Item(68, 8).Value = 2099.1428  # H68: 2169.2 -> 2099.1428
$ws.Cells.Item(68, 9).Value = 2159  # I68: 2224 -> 2159
$ws.Cells.Item(68, 10).Value = 1949.5  # J68: 1950 -> 1949.5
$ws.Cells.Item(68, 11).Value = 2159  # K68: 2224 -> 2159
$ws.Cells.Item(68, 12).Value = 1949.5  # L68: 1950 -> 1949.5
$ws.Cells.Item(68, 13).Value = -1410  # M68: -1475 -> -1410
$ws.Cells.Item(68, 14).Value = -3447.5  # N68: -3448 -> -3447.5
$ws.Cells.Item(71, 8).Value = 2099.1428  # H71: 2169.2 -> 2099.1428
$ws.Cells.Item(71, 9).Value = 2159  # I71: 2224 -> 2159
$ws.Cells.Item(71, 10).Value = 1949.5  # J71: 1950 -> 1949.5
$ws.Cells.Item(71, 11).Value = 10795  # K71: 11120 -> 10795
$ws.Cells.Item(71, 12).Value = 9747.5  # L71: 9750 -> 9747.5
$ws.Cells.Item(71, 13).Value = -7051  # M71: -7376 -> -7051
$ws.Cells.Item(71, 14).Value = -17235.5  # N71: -17238 -> -17235.5
$ws.Cells.Item(80, 8).Value = 38331.332  # H80: 39997 -> 38331.332
$ws.Cells.Item(80, 9).Value = 35000  # I80: 0 -> 35000
$ws.Cells.Item(80, 11).Value = 35000  # K80: 0 -> 35000
$ws.Cells.Item(80, 13).Value = -33877  # M80: None -> -33877
$ws.Cells.Item(83, 8).Value = 38331.332  # H83: 39997 -> 38331.332
$ws.Cells.Item(83, 9).Value = 35000  # I83: 0 -> 35000
$ws.Cells.Item(83, 11).Value = 105000  # K83: 0 -> 105000
$ws.Cells.Item(83, 13).Value = -99384  # M83: None -> -99384
$ws.Cells.Item(113, 8).Value = 7418.8335  # H113: 7423 -> 7418.8335
$ws.Cells.Item(113, 9).Value = 7422.8  # I113: 7423 -> 7422.8
$ws.Cells.Item(113, 10).Value = 7399  # J113: 0 -> 7399
$ws.Cells.Item(113, 11).Value = 7422.8  # K113: 7423 -> 7422.8
$ws.Cells.Item(113, 12).Value = 7399  # L113: 0 -> 7399
$ws.Cells.Item(113, 13).Value = -5252.8  # M113: -5253 -> -5252.8
$ws.Cells.Item(113, 14).Value = -11739  # N113: None -> -11739
$ws.Cells.Item(122, 8).Value = 4499  # H122: 4998.75 -> 4499
$ws.Cells.Item(122, 9).Value = 3002  # I122: 3504 -> 3002
$ws.Cells.Item(122, 11).Value = 9006  # K122: 10512 -> 9006
$ws.Cells.Item(122, 13).Value = -6556  # M122: -8062 -> -6556
$ws.Cells.Item(126, 8).Value = 3947  # H126: 0 -> 3947
$ws.Cells.Item(126, 9).Value = 3904  # I126: 0 -> 3904
$ws.Cells.Item(126, 10).Value = 3990  # J126: 0 -> 3990
$ws.Cells.Item(126, 11).Value = 11712  # K126: 0 -> 11712
$ws.Cells.Item(126, 12).Value = 11970  # L126: 0 -> 11970
$ws.Cells.Item(126, 13).Value = -9242  # M126: None -> -9242
$ws.Cells.Item(126, 14).Value = -16910  # N126: None -> -16910
$ws.Cells.Item(132, 8).Value = 9700.299999999999  # H132: 9799.799999999999 -> 9700.299999999999
$ws.Cells.Item(132, 9).Value = 7947.8423  # I132: 8138.6113 -> 7947.8423
$ws.Cells.Item(132, 10).Value = 12727.272  # J132: 12291.583 -> 12727.272
$ws.Cells.Item(132, 11).Value = 23843.5269  # K132: 24415.8339 -> 23843.5269
$ws.Cells.Item(132, 12).Value = 38181.81600000001  # L132: 36874.749 -> 38181.81600000001
$ws.Cells.Item(132, 13).Value = -21313.5269  # M132: -21885.8339 -> -21313.5269
$ws.Cells.Item(132, 14).Value = -43241.81600000001  # N132: -41934.749 -> -43241.81600000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1601  # H113: 880.3 -> 1601
$ws.Cells.Item(113, 9).Value = 2302  # I113: 900.4286 -> 2302
$ws.Cells.Item(113, 10).Value = 900  # J113: 833.3333 -> 900
$ws.Cells.Item(113, 11).Value = 6906  # K113: 2701.2858 -> 6906
$ws.Cells.Item(113, 12).Value = 2700  # L113: 2499.9999 -> 2700
$ws.Cells.Item(113, 13).Value = -4736  # M113: -531.2857999999997 -> -4736
$ws.Cells.Item(113, 14).Value = -7040  # N113: -6839.9999 -> -7040
$ws.Cells.Item(126, 8).Value = 2050.75  # H126: 2375 -> 2050.75
$ws.Cells.Item(126, 9).Value = 2050.75  # I126: 2375 -> 2050.75
$ws.Cells.Item(126, 11).Value = 6152.25  # K126: 7125 -> 6152.25
$ws.Cells.Item(126, 13).Value = -3682.25  # M126: -4655 -> -3682.25
$ws.Cells.Item(132, 8).Value = 2005.9286  # H132: 1775.6923 -> 2005.9286
$ws.Cells.Item(132, 10).Value = 4999.5  # J132: 5000 -> 4999.5
$ws.Cells.Item(132, 12).Value = 14998.5  # L132: 15000 -> 14998.5
$ws.Cells.Item(132, 14).Value = -20058.5  # N132: -20060 -> -20058.5
